$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.362.44'
$ws.Range("E2").Value = '  +0.68%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.932.63'
$ws.Range("E3").Value = '  +0.61%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.05'
$ws.Range("E5").Value = '  +0.37%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.77'
$ws.Range("E6").Value = '  -1.49%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("E8").Value = '  -1.21%  '

$ws.Range("E9").Value = '  +0.89%  '

$ws.Range("E10").Value = '  -1.42%  '

$ws.Range("E11").Value = '  -0.49%  '

$ws.Range("E12").Value = '  -0.89%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '33.25'
$ws.Range("E13").Value = '  -0.85%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.418.30'
$ws.Range("E15").Value = '  +0.63%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '61.402.85'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.933.41'
$ws.Range("E17").Value = '  +0.58%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.64'
$ws.Range("E18").Value = '  -0.80%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '433.68'
$ws.Range("E19").Value = '  +0.72%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.56'
$ws.Range("E20").Value = '  +1.48%  '

$ws.Range("E21").Value = '  -0.93%  '

$ws.Range("E22").Value = '  +0.22%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '81.42'
$ws.Range("E23").Value = '  -0.26%  '

$ws.Range("E24").Value = '  -0.27%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.17'
$ws.Range("E25").Value = '  -1.39%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.75'
$ws.Range("E26").Value = '  -1.30%  '

$ws.Range("E28").Value = '  -3.82%  '

$ws.Range("E29").Value = '  -0.48%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.89'
$ws.Range("E30").Value = '  -2.20%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '26.76'
$ws.Range("E31").Value = '  +0.93%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.109'
$ws.Range("E32").Value = '  +1.30%  '

$ws.Range("E33").Value = '  +0.09%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0₃0874'
$ws.Range("E34").Value = '  +2.54%  '

$ws.Range("E35").Value = '  -0.11%  '

$ws.Range("E36").Value = '  -0.25%  '

$ws.Range("E37").Value = '  -1.48%  '

$ws.Range("E38").Value = '  +0.24%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.123'
$ws.Range("E39").Value = '  -0.50%  '

$ws.Range("E40").Value = '  -0.56%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '42.07'
$ws.Range("E41").Value = '  +5.35%  '

$ws.Range("E42").Value = '  -1.98%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0345'
$ws.Range("E43").Value = '  +0.06%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.695.93'
$ws.Range("E44").Value = '  -0.33%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '364.78'
$ws.Range("E45").Value = '  -2.81%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '133.49'
$ws.Range("E46").Value = '  +0.92%  '

$ws.Range("E47").Value = '  +0.11%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.55'
$ws.Range("E48").Value = '  -1.28%  '

$ws.Range("E49").Value = '  -1.32%  '

$ws.Range("E50").Value = '  -0.54%  '

$ws.Range("E51").Value = '  +0.26%  '
